$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.004953627445374112
$ws.Range("J2").Value = 0.004953627445374111
$ws.Range("M2").Value = 0.2753413333333334
$ws.Range("N2").Value = 0.8260240000000001
$ws.Range("O2").Value = 0.006630378892106956
$ws.Range("P2").Value = 0.006630378892106955
$ws.Range("Q2").Value = 0.0009779206355555556
$ws.Range("R2").Value = 0.008801285720000001
$ws.Range("S2").Value = 0.00003284442685317021
$ws.Range("T2").Value = 0.00003284442685317021
$ws.Range("I3").Value = 0.004953627445374112
$ws.Range("J3").Value = 0.004953627445374111
$ws.Range("O3").Value = 0.03952244389885164
$ws.Range("P3").Value = 0.03952244389885164
$ws.Range("S3").Value = 0.0001957794628056101
$ws.Range("T3").Value = 0.0001957794628056101
$ws.Range("I4").Value = 0.004953627445374112
$ws.Range("J4").Value = 0.004953627445374111
$ws.Range("M4").Value = 23.78768866666667
$ws.Range("N4").Value = 71.363066
$ws.Range("O4").Value = 0.5728213302306416
$ws.Range("P4").Value = 0.5728213302306416
$ws.Range("Q4").Value = 0.08448594091444445
$ws.Range("R4").Value = 0.76037346823
$ws.Range("S4").Value = 0.002837543462726214
$ws.Range("T4").Value = 0.002837543462726213
$ws.Range("I5").Value = 0.004953627445374112
$ws.Range("J5").Value = 0.004953627445374111
$ws.Range("M5").Value = 0.5982033333333333
$ws.Range("N5").Value = 1.79461
$ws.Range("O5").Value = 0.01440508298011203
$ws.Range("P5").Value = 0.01440508298011203
$ws.Range("Q5").Value = 0.002124618838888889
$ws.Range("R5").Value = 0.01912156955
$ws.Range("S5").Value = 0.00007135741440317447
$ws.Range("T5").Value = 0.00007135741440317446
$ws.Range("I6").Value = 0.004953627445374112
$ws.Range("J6").Value = 0.004953627445374111
$ws.Range("M6").Value = 15.22474833333333
$ws.Range("N6").Value = 45.674245
$ws.Range("O6").Value = 0.3666207639982877
$ws.Range("P6").Value = 0.3666207639982877
$ws.Range("Q6").Value = 0.05407323116388888
$ws.Range("R6").Value = 0.4866590804749999
$ws.Range("S6").Value = 0.001816102678585943
$ws.Range("T6").Value = 0.001816102678585943
$ws.Range("I7").Value = 0.0778289024983856
$ws.Range("J7").Value = 0.07782890249838559
$ws.Range("M7").Value = 0.2753413333333334
$ws.Range("N7").Value = 0.8260240000000001
$ws.Range("O7").Value = 0.006630378892106956
$ws.Range("P7").Value = 0.006630378892106955
$ws.Range("Q7").Value = 0.01536459708266667
$ws.Range("R7").Value = 0.138281373744
$ws.Range("S7").Value = 0.0005160351123211463
$ws.Range("T7").Value = 0.0005160351123211461
$ws.Range("I8").Value = 0.0778289024983856
$ws.Range("J8").Value = 0.07782890249838559
$ws.Range("O8").Value = 0.03952244389885164
$ws.Range("P8").Value = 0.03952244389885164
$ws.Range("S8").Value = 0.003075988432701639
$ws.Range("T8").Value = 0.003075988432701639
$ws.Range("I9").Value = 0.0778289024983856
$ws.Range("J9").Value = 0.07782890249838559
$ws.Range("M9").Value = 23.78768866666667
$ws.Range("N9").Value = 71.363066
$ws.Range("O9").Value = 0.5728213302306416
$ws.Range("P9").Value = 0.5728213302306416
$ws.Range("Q9").Value = 1.327400602977333
$ws.Range("R9").Value = 11.946605426796
$ws.Range("S9").Value = 0.04458205545951615
$ws.Range("T9").Value = 0.04458205545951614
$ws.Range("I10").Value = 0.0778289024983856
$ws.Range("J10").Value = 0.07782890249838559
$ws.Range("M10").Value = 0.5982033333333333
$ws.Range("N10").Value = 1.79461
$ws.Range("O10").Value = 0.01440508298011203
$ws.Range("P10").Value = 0.01440508298011203
$ws.Range("Q10").Value = 0.03338094240666666
$ws.Range("R10").Value = 0.30042848166
$ws.Range("S10").Value = 0.001121131798740293
$ws.Range("T10").Value = 0.001121131798740293
$ws.Range("I11").Value = 0.0778289024983856
$ws.Range("J11").Value = 0.07782890249838559
$ws.Range("M11").Value = 15.22474833333333
$ws.Range("N11").Value = 45.674245
$ws.Range("O11").Value = 0.3666207639982877
$ws.Range("P11").Value = 0.3666207639982877
$ws.Range("Q11").Value = 0.8495714064966666
$ws.Range("R11").Value = 7.64614265847
$ws.Range("S11").Value = 0.02853369169510638
$ws.Range("T11").Value = 0.02853369169510636
$ws.Range("G12").Value = 0.509521
$ws.Range("H12").Value = 1.528563
$ws.Range("I12").Value = 0.7106458591068409
$ws.Range("J12").Value = 0.7106458591068406
$ws.Range("M12").Value = 0.2753413333333334
$ws.Range("N12").Value = 0.8260240000000001
$ws.Range("O12").Value = 0.006630378892106956
$ws.Range("P12").Value = 0.006630378892106955
$ws.Range("Q12").Value = 0.1402921915013333
$ws.Range("R12").Value = 1.262629723512
$ws.Range("S12").Value = 0.004711851303985212
$ws.Range("T12").Value = 0.00471185130398521
$ws.Range("G13").Value = 0.509521
$ws.Range("H13").Value = 1.528563
$ws.Range("I13").Value = 0.7106458591068409
$ws.Range("J13").Value = 0.7106458591068406
$ws.Range("O13").Value = 0.03952244389885164
$ws.Range("P13").Value = 0.03952244389885164
$ws.Range("Q13").Value = 0.8362554174179999
$ws.Range("R13").Value = 7.526298756761999
$ws.Range("S13").Value = 0.02808646109850135
$ws.Range("T13").Value = 0.02808646109850134
$ws.Range("G14").Value = 0.509521
$ws.Range("H14").Value = 1.528563
$ws.Range("I14").Value = 0.7106458591068409
$ws.Range("J14").Value = 0.7106458591068406
$ws.Range("M14").Value = 23.78768866666667
$ws.Range("N14").Value = 71.363066
$ws.Range("O14").Value = 0.5728213302306416
$ws.Range("P14").Value = 0.5728213302306416
$ws.Range("Q14").Value = 12.12032691712867
$ws.Range("R14").Value = 109.082942254158
$ws.Range("S14").Value = 0.4070731063364777
$ws.Range("T14").Value = 0.4070731063364776
$ws.Range("G15").Value = 0.509521
$ws.Range("H15").Value = 1.528563
$ws.Range("I15").Value = 0.7106458591068409
$ws.Range("J15").Value = 0.7106458591068406
$ws.Range("M15").Value = 0.5982033333333333
$ws.Range("N15").Value = 1.79461
$ws.Range("O15").Value = 0.01440508298011203
$ws.Range("P15").Value = 0.01440508298011203
$ws.Range("Q15").Value = 0.3047971606033333
$ws.Range("R15").Value = 2.74317444543
$ws.Range("S15").Value = 0.01023691256990705
$ws.Range("T15").Value = 0.01023691256990704
$ws.Range("G16").Value = 0.509521
$ws.Range("H16").Value = 1.528563
$ws.Range("I16").Value = 0.7106458591068409
$ws.Range("J16").Value = 0.7106458591068406
$ws.Range("M16").Value = 15.22474833333333
$ws.Range("N16").Value = 45.674245
$ws.Range("O16").Value = 0.3666207639982877
$ws.Range("P16").Value = 0.3666207639982877
$ws.Range("Q16").Value = 7.757328995548333
$ws.Range("R16").Value = 69.815960959935
$ws.Range("S16").Value = 0.2605375277979696
$ws.Range("T16").Value = 0.2605375277979694
$ws.Range("G17").Value = 0.06794533333333334
$ws.Range("H17").Value = 0.203836
$ws.Range("I17").Value = 0.09476561275976328
$ws.Range("J17").Value = 0.09476561275976325
$ws.Range("M17").Value = 0.2753413333333334
$ws.Range("N17").Value = 0.8260240000000001
$ws.Range("O17").Value = 0.006630378892106956
$ws.Range("P17").Value = 0.006630378892106955
$ws.Range("Q17").Value = 0.01870815867377778
$ws.Range("R17").Value = 0.168373428064
$ws.Range("S17").Value = 0.0006283319185399161
$ws.Range("T17").Value = 0.0006283319185399158
$ws.Range("G18").Value = 0.06794533333333334
$ws.Range("H18").Value = 0.203836
$ws.Range("I18").Value = 0.09476561275976328
$ws.Range("J18").Value = 0.09476561275976325
$ws.Range("O18").Value = 0.03952244389885164
$ws.Range("P18").Value = 0.03952244389885164
$ws.Range("Q18").Value = 0.111515821896
$ws.Range("R18").Value = 1.003642397064
$ws.Range("S18").Value = 0.003745368613838044
$ws.Range("T18").Value = 0.003745368613838043
$ws.Range("G19").Value = 0.06794533333333334
$ws.Range("H19").Value = 0.203836
$ws.Range("I19").Value = 0.09476561275976328
$ws.Range("J19").Value = 0.09476561275976325
$ws.Range("M19").Value = 23.78768866666667
$ws.Range("N19").Value = 71.363066
$ws.Range("O19").Value = 0.5728213302306416
$ws.Range("P19").Value = 0.5728213302306416
$ws.Range("Q19").Value = 1.616262435686223
$ws.Range("R19").Value = 14.546361921176
$ws.Range("S19").Value = 0.05428376436116947
$ws.Range("T19").Value = 0.05428376436116945
$ws.Range("G20").Value = 0.06794533333333334
$ws.Range("H20").Value = 0.203836
$ws.Range("I20").Value = 0.09476561275976328
$ws.Range("J20").Value = 0.09476561275976325
$ws.Range("M20").Value = 0.5982033333333333
$ws.Range("N20").Value = 1.79461
$ws.Range("O20").Value = 0.01440508298011203
$ws.Range("P20").Value = 0.01440508298011203
$ws.Range("Q20").Value = 0.04064512488444445
$ws.Range("R20").Value = 0.3658061239600001
$ws.Range("S20").Value = 0.001365106515465554
$ws.Range("T20").Value = 0.001365106515465553
$ws.Range("G21").Value = 0.06794533333333334
$ws.Range("H21").Value = 0.203836
$ws.Range("I21").Value = 0.09476561275976328
$ws.Range("J21").Value = 0.09476561275976325
$ws.Range("M21").Value = 15.22474833333333
$ws.Range("N21").Value = 45.674245
$ws.Range("O21").Value = 0.3666207639982877
$ws.Range("P21").Value = 0.3666207639982877
$ws.Range("Q21").Value = 1.034450600424444
$ws.Range("R21").Value = 9.31005540382
$ws.Range("S21").Value = 0.0347430413507503
$ws.Range("T21").Value = 0.03474304135075029
$ws.Range("G22").Value = 0.080163
$ws.Range("H22").Value = 0.240489
$ws.Range("I22").Value = 0.1118059981896363
$ws.Range("J22").Value = 0.1118059981896363
$ws.Range("M22").Value = 0.2753413333333334
$ws.Range("N22").Value = 0.8260240000000001
$ws.Range("O22").Value = 0.006630378892106956
$ws.Range("P22").Value = 0.006630378892106955
$ws.Range("Q22").Value = 0.022072187304
$ws.Range("R22").Value = 0.198649685736
$ws.Range("S22").Value = 0.0007413161304075131
$ws.Range("T22").Value = 0.0007413161304075128
$ws.Range("G23").Value = 0.080163
$ws.Range("H23").Value = 0.240489
$ws.Range("I23").Value = 0.1118059981896363
$ws.Range("J23").Value = 0.1118059981896363
$ws.Range("O23").Value = 0.03952244389885164
$ws.Range("P23").Value = 0.03952244389885164
$ws.Range("Q23").Value = 0.131568165054
$ws.Range("R23").Value = 1.184113485486
$ws.Range("S23").Value = 0.00441884629100501
$ws.Range("T23").Value = 0.004418846291005008
$ws.Range("G24").Value = 0.080163
$ws.Range("H24").Value = 0.240489
$ws.Range("I24").Value = 0.1118059981896363
$ws.Range("J24").Value = 0.1118059981896363
$ws.Range("M24").Value = 23.78768866666667
$ws.Range("N24").Value = 71.363066
$ws.Range("O24").Value = 0.5728213302306416
$ws.Range("P24").Value = 0.5728213302306416
$ws.Range("Q24").Value = 1.906892486586
$ws.Range("R24").Value = 17.162032379274
$ws.Range("S24").Value = 0.06404486061075218
$ws.Range("T24").Value = 0.06404486061075215
$ws.Range("G25").Value = 0.080163
$ws.Range("H25").Value = 0.240489
$ws.Range("I25").Value = 0.1118059981896363
$ws.Range("J25").Value = 0.1118059981896363
$ws.Range("M25").Value = 0.5982033333333333
$ws.Range("N25").Value = 1.79461
$ws.Range("O25").Value = 0.01440508298011203
$ws.Range("P25").Value = 0.01440508298011203
$ws.Range("Q25").Value = 0.04795377381
$ws.Range("R25").Value = 0.43158396429
$ws.Range("S25").Value = 0.001610574681595967
$ws.Range("T25").Value = 0.001610574681595966
$ws.Range("G26").Value = 0.080163
$ws.Range("H26").Value = 0.240489
$ws.Range("I26").Value = 0.1118059981896363
$ws.Range("J26").Value = 0.1118059981896363
$ws.Range("M26").Value = 15.22474833333333
$ws.Range("N26").Value = 45.674245
$ws.Range("O26").Value = 0.3666207639982877
$ws.Range("P26").Value = 0.3666207639982877
$ws.Range("Q26").Value = 1.220461500645
$ws.Range("R26").Value = 10.984153505805
$ws.Range("S26").Value = 0.04099040047587564
$ws.Range("T26").Value = 0.04099040047587562
